# FedExShipments_NewRateChng.xlsx - "Changes of 28th july 2022"
#
# New FedEx tracking numbers were issued for the shipments on rows 2-26
# (column P), and the rate-check result for row 20 flipped from PASS to
# FAIL (with updated expected/actual rate figures) while row 24 flipped
# from FAIL to PASS.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracking numbers (FedEx "ShipmentTracking" column) for rows 2-26.
$trackingNumbers = @{
    2  = "320018767744"
    3  = "320018767814"
    4  = "320018767847"
    5  = "320018767869"
    6  = "320018767906"
    7  = "320018767928"
    8  = "320018768155"
    9  = "320018768177"
    10 = "320018768203"
    11 = "320018768225"
    12 = "320018768269"
    13 = "320018768280"
    14 = "320018768317"
    15 = "320018768339"
    16 = "320018768361"
    17 = "320018768383"
    18 = "320018768420"
    19 = "320018768442"
    20 = "320018768475"
    21 = "320018768497"
    22 = "320018768523"
    23 = "320018768534"
    24 = "320018768545"
    25 = "320018768556"
    26 = "320018768567"
}

# These tracking numbers are all-digit strings; Excel would normally infer
# a number when assigning through .Value, so the destination cells are
# temporarily marked as Text before the write (and restored afterwards)
# to keep them stored as plain shared-string text, matching the rest of
# the column.
$trackRange = $ws.Range("P2:P26")
$trackRange.NumberFormat = "@"
foreach ($row in $trackingNumbers.Keys) {
    $ws.Range("P$row").Value = $trackingNumbers[$row]
}
$trackRange.ClearFormats()

# Row 20's rate check now fails: the actual rate rose to $109.19 against
# the shipment's expected rate, and the per-pound surcharge became $5.94.
$dollarCells = @{
    "Q20" = "`$109.19"
    "U20" = "`$5.94"
    "W20" = "`$109.19"
}
foreach ($ref in $dollarCells.Keys) {
    $ws.Range($ref).NumberFormat = "@"
}
foreach ($ref in $dollarCells.Keys) {
    $ws.Range($ref).Value = $dollarCells[$ref]
}
foreach ($ref in $dollarCells.Keys) {
    $ws.Range($ref).ClearFormats()
}

$ws.Range("R20").Value = "FAIL"

# Row 24 now passes the rate check.
$ws.Range("R24").Value = "PASS"
